# Rebuild the invoices sheet:
#  - relabel the header row (now Greek display labels, plus a new "Total" column)
#  - replace the data rows with the new values (every value now arrives as
#    plain text from the back end, including the numeric-looking ones)
#  - drop the old 4th data row entirely (3 data rows -> 2 data rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old 4th sheet row (the extra duplicate invoice line) ---
$ws.Rows.Item(4).Delete()

# --- New column K: seed it from column J so it inherits the same header
#     (bold/centered/bordered) style as the rest of row 1. ---
$ws.Cells.Item(1,10).Copy($ws.Cells.Item(1,11))

# --- Header row (row 1) — keep the existing style, just relabel the text ---
$headers = @("MARK","ΑΦΜ","Επωνυμία","Σειρά","Αριθμός","Ημερομηνία","Είδος","ΦΠΑ_ΚΑΤΗΓΟΡΙΑ","Καθαρή Αξία","ΦΠΑ","Σύνολο")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value2 = $headers[$col - 1]
}

# --- Data rows -------------------------------------------------------------
# Every incoming value is plain text (no more back-end JSON typing), even the
# numeric-looking ones (amounts, counts, dates). Force text storage with a
# temporary Text number format before writing the value, mirroring how
# Excel keeps a value as text instead of auto-converting it to a number or
# date, then drop the number-format override again so no stray per-cell
# style sticks around once the value is in place.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = "Normal"
}

$row2 = @("400008186167654", "", "", "8Μ0ΤΔΑ", "", "2025-01-03", "", "1", "75.98", "18.24", "94.22")
$row3 = @("400008186167654", "", "", "8Μ0ΤΔΑ", "", "2025-01-03", "", "1", "75.98", "18.24", "94.22")

for ($col = 1; $col -le $row2.Length; $col++) {
    Set-TextValue $ws.Cells.Item(2, $col) $row2[$col - 1]
}
for ($col = 1; $col -le $row3.Length; $col++) {
    Set-TextValue $ws.Cells.Item(3, $col) $row3[$col - 1]
}

Write-Output ("UsedRange=" + $ws.UsedRange.Address())
